# Applies the "automatic update" commit:
#  1) Bumps the "Förändrad" (changed) date in column C for every existing
#     data row (2..300) from 2023-09-13 (45182) to 2023-09-15 (45184).
#  2) Gives row 300 an explicit custom row height (ht="15" customHeight="1"),
#     matching what Excel stamps on a row once it has been "touched" again.
#  3) Appends a brand-new record as row 301 (case "A 42986-2023").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") for all existing data rows.
$ws.Range("C2:C300").Value = 45184

# 2) Stamp row 300 with an explicit row height (no visual change, 15pt is
#    already the sheet default, but it becomes an explicit/custom height).
$ws.Rows.Item(300).RowHeight = 15

# 3) Append the new row (301) with its data.
$newRow = 301

$ws.Cells.Item($newRow, 1).Value = "A 42986-2023"

$ws.Cells.Item($newRow, 2).Value = 45182
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45184
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($newRow, 5).Value = "BORLÄNGE"

# Column F (Markägare) intentionally left blank, as with every other row.

$ws.Cells.Item($newRow, 7).Value = 0.8

$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Column R (Artnamn) keeps the same wrap-text styling used for every other
# row's (empty) species-name cell.
$ws.Cells.Item($newRow, 18).WrapText = $true
